# Add New TC for Sim ATM
# - Adds a new value "run" in cell A2 (new shared string)
# - Moves the active selection from G2 to A2

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "run"
$ws.Range("A2").Select()
